# added harvard case classification
# Recompute the "_old" (pre-harvard-case-classification) stats columns plus the
# average_doctor / average_doctor_old pair, and flip which column each of
# those two headers/labels refers to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers: the "average_doctor" / "average_doctor_old" labels swap columns.
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4 (stats_for_precision)
$ws.Range("E4").Value  = 0.54
$ws.Range("F4").Value  = 0.066
$ws.Range("G4").Value  = 0.258
$ws.Range("N4").Value  = 0.367
$ws.Range("W4").Value  = 0.2
$ws.Range("X4").Value  = 0.071
$ws.Range("Y4").Value  = 0.267
$ws.Range("AI4").Value = 0.12
$ws.Range("AJ4").Value = 0.026
$ws.Range("AK4").Value = 0.16
$ws.Range("BA4").Value = 1.35
$ws.Range("BB4").Value = 0.14
$ws.Range("BC4").Value = 0.374
$ws.Range("BG4").Value = 0.5
$ws.Range("BH4").Value = 0.2
$ws.Range("BI4").Value = 0.447
$ws.Range("BP4").Value = 0.45
$ws.Range("BQ4").Value = 0.551

# Row 5 (stats_for_recall)
$ws.Range("E5").Value  = 0.517
$ws.Range("F5").Value  = 0.029
$ws.Range("G5").Value  = 0.17
$ws.Range("N5").Value  = 0.75
$ws.Range("O5").Value  = 0.061
$ws.Range("P5").Value  = 0.247
$ws.Range("W5").Value  = 0.167
$ws.Range("X5").Value  = 0.044
$ws.Range("Y5").Value  = 0.211
$ws.Range("AI5").Value = 0.167
$ws.Range("AJ5").Value = 0.044
$ws.Range("AK5").Value = 0.211
$ws.Range("BA5").Value = 0.95
$ws.Range("BB5").Value = 0.051
$ws.Range("BC5").Value = 0.226
$ws.Range("BG5").Value = 0.333
$ws.Range("BH5").Value = 0.078
$ws.Range("BI5").Value = 0.279
$ws.Range("BP5").Value = 0.317
$ws.Range("BQ5").Value = 0.366

# Row 6 (stats_for_f1-score)
$ws.Range("E6").Value  = 0.528
$ws.Range("N6").Value  = 0.493
$ws.Range("W6").Value  = 0.182
$ws.Range("AI6").Value = 0.14
$ws.Range("BA6").Value = 1.115
$ws.Range("BG6").Value = 0.4
$ws.Range("BP6").Value = 0.372
$ws.Range("BQ6").Value = 0.438

# Row 7 (stats_for_f2-score)
$ws.Range("E7").Value  = 0.521
$ws.Range("N7").Value  = 0.62
$ws.Range("W7").Value  = 0.173
$ws.Range("AI7").Value = 0.155
$ws.Range("BA7").Value = 1.009
$ws.Range("BG7").Value = 0.357
$ws.Range("BP7").Value = 0.336
$ws.Range("BQ7").Value = 0.392

# Row 8 (stats_for_NDCG)
$ws.Range("E8").Value  = 0.489
$ws.Range("F8").Value  = 0.055
$ws.Range("G8").Value  = 0.234
$ws.Range("N8").Value  = 0.917
$ws.Range("O8").Value  = 0.009
$ws.Range("P8").Value  = 0.095
$ws.Range("W8").Value  = 0.221
$ws.Range("X8").Value  = 0.088
$ws.Range("Y8").Value  = 0.297
$ws.Range("AI8").Value = 0.073
$ws.Range("AJ8").Value = 0.009
$ws.Range("AK8").Value = 0.093
$ws.Range("BA8").Value = 1.339
$ws.Range("BB8").Value = 0.149
$ws.Range("BC8").Value = 0.386
$ws.Range("BG8").Value = 0.403
$ws.Range("BH8").Value = 0.158
$ws.Range("BI8").Value = 0.397
$ws.Range("BP8").Value = 0.446
$ws.Range("BQ8").Value = 0.503

# Row 9 (stats_for_M1) - average_doctor / average_doctor_old swap only
$ws.Range("BP9").Value  = 0.467
$ws.Range("BQ9").Value  = 0.49

# Row 10 (stats_for_M3)
$ws.Range("BP10").Value = 0.533
$ws.Range("BQ10").Value = 0.647

# Row 11 (stats_for_M5)
$ws.Range("BP11").Value = 0.533
$ws.Range("BQ11").Value = 0.647

# Row 12 (stats_for_position)
$ws.Range("BP12").Value = 1.083
$ws.Range("BQ12").Value = 1.364

# Row 13 (stats_for_length (x of gs))
$ws.Range("BP13").Value = 0.778
$ws.Range("BQ13").Value = 0.753
